# Edit: rename headers, add new "PO Forecast" sheet with forecast data.

$wb = $excel.ActiveWorkbook

# --- 1) Rename header cells on the two existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row, bold + centered + top-aligned + thin box border (matches the
# other sheets' header style).
$header = $wsForecast.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Data rows
$forecastRows = @(
    @(2, 45326.99999999999, 27, -9.638268051100155, 61.15240943972568),
    @(3, 45354.99999999999, 26, -10.20756827551194, 61.14598864586915),
    @(4, 45368.99999999999, 26, -8.624825912384008, 62.07935835853147),
    @(5, 45375.99999999999, 26, -9.672291485610213, 59.74933779034938),
    @(6, 45389.99999999999, 25, -8.844192257393384, 61.46274189131418),
    @(7, 45410.99999999999, 25, -11.77235668766766, 63.76462909841526),
    @(8, 45424.99999999999, 25, -10.1163786217858, 61.26860499263336),
    @(9, 45431.99999999999, 25, -10.99203433545378, 59.67096506437193),
    @(10, 45438.99999999999, 24, -12.89810431913743, 59.59058756818318),
    @(11, 45445.99999999999, 24, -9.610195313638316, 60.85947221128671),
    @(12, 45452.99999999999, 24, -14.09680687036047, 60.43481242367945),
    @(13, 45459.99999999999, 24, -13.295710847413, 60.19208906112828),
    @(14, 45473.99999999999, 24, -13.0079293889472, 59.96500888550218),
    @(15, 45480.99999999999, 24, -10.90121448246101, 61.03462750014179),
    @(16, 45529.99999999999, 23, -13.09560967576719, 59.98664345099824),
    @(17, 45536.99999999999, 23, -11.80268324854338, 59.30640131732664),
    @(18, 45543.99999999999, 22, -14.02513056969536, 60.9738988188676),
    @(19, 45550.99999999999, 22, -13.86759305366119, 57.51486018704174),
    @(20, 45578.99999999999, 22, -16.60052909824957, 55.81537898984726),
    @(21, 45585.99999999999, 22, -14.69181573517072, 61.66936717342635),
    @(22, 45620.99999999999, 21, -17.76524775076977, 57.61205414919857),
    @(23, 45627.99999999999, 21, -16.53515966867629, 58.81049910758249),
    @(24, 45634.99999999999, 21, -14.71170869458271, 54.51453864075712),
    @(25, 45641.99999999999, 20, -16.04429444589852, 58.26248093589494),
    @(26, 45648.99999999999, 20, -19.55095724937767, 56.50793095517382),
    @(27, 45655.99999999999, 20, -17.04097639001375, 57.28556037254592),
    @(28, 45662.99999999999, 20, -13.40586321213311, 58.37039596217429),
    @(29, 45669.99999999999, 20, -15.83587694146827, 57.33500337332152),
    @(30, 45676.99999999999, 20, -14.7354208707975, 56.34774101562561),
    @(31, 45683.99999999999, 20, -17.6947861897221, 56.86128231670242),
    @(32, 45690.99999999999, 19, -17.01280849885342, 55.21898307256797),
    @(33, 45697.99999999999, 19, -14.92947158741028, 55.55357017369288)
)

foreach ($row in $forecastRows) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}
